$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster changes from FAPs -> ECs, plus refreshed metrics ---
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.8330250000000001
$ws.Cells.Item(2, 8).Value = 2.499075
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.4056473333333333
$ws.Cells.Item(2, 14).Value = 1.216942
$ws.Cells.Item(2, 15).Value = 0.2120329281749088
$ws.Cells.Item(2, 16).Value = 0.2120329281749088
$ws.Cells.Item(2, 17).Value = 0.33791436985
$ws.Cells.Item(2, 18).Value = 3.041229328650001
$ws.Cells.Item(2, 19).Value = 0.2120329281749088
$ws.Cells.Item(2, 20).Value = 0.2120329281749088

# --- Row 3: Target cluster changes from sCs -> FAPs, plus refreshed metrics ---
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.8330250000000001
$ws.Cells.Item(3, 8).Value = 2.499075
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.086432
$ws.Cells.Item(3, 14).Value = 3.259296
$ws.Cells.Item(3, 15).Value = 0.5678808642225905
$ws.Cells.Item(3, 16).Value = 0.5678808642225904
$ws.Cells.Item(3, 17).Value = 0.9050250168000001
$ws.Cells.Item(3, 18).Value = 8.145225151200002
$ws.Cells.Item(3, 19).Value = 0.5678808642225905
$ws.Cells.Item(3, 20).Value = 0.5678808642225904

# --- Row 4 (new): Sending cluster FAPs / Ligand Wnt2 / Receptor Fzd9 / Target cluster sCs ---
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Wnt2"
$ws.Cells.Item(4, 3).Value = "Fzd9"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.8330250000000001
$ws.Cells.Item(4, 8).Value = 2.499075
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4210543333333334
$ws.Cells.Item(4, 14).Value = 1.263163
$ws.Cells.Item(4, 15).Value = 0.2200862076025007
$ws.Cells.Item(4, 16).Value = 0.2200862076025006
$ws.Cells.Item(4, 17).Value = 0.3507487860250001
$ws.Cells.Item(4, 18).Value = 3.156739074225001
$ws.Cells.Item(4, 19).Value = 0.2200862076025007
$ws.Cells.Item(4, 20).Value = 0.2200862076025006
